# Updated cryptos list on Fri May  3 03:56:22 UTC 2024 with GitHub Actions
#
# Applies the per-row price / volume(1h) refresh to the crypto price sheet,
# plus the two swapped rows (NEARProtocol <-> EthereumClassic, rows 31/32).
#
# Values are written as literal text (matching the source inlineStr cells):
# NumberFormat is forced to "@" (Text) before the write so Excel does not
# silently coerce strings like "566.30" or "5.31" into numeric values
# (which would drop trailing zeros / change precision), and the cell's
# original Style is restored afterwards so no stray style index is left
# behind on cells that were previously unstyled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetCellText {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = $origStyle
}

SetCellText 'D2' '59.899.93'
SetCellText 'E2' '  +4.05%  '
SetCellText 'D3' '3.020.75'
SetCellText 'E3' '  +3.10%  '
SetCellText 'E4' '  +0.18%  '
SetCellText 'D5' '566.30'
SetCellText 'E5' '  +3.03%  '
SetCellText 'D6' '141.23'
SetCellText 'E6' '  +8.40%  '
SetCellText 'E7' '  -0.03%  '
SetCellText 'E8' '  +2.12%  '
SetCellText 'D9' '3.009.61'
SetCellText 'E9' '  +2.79%  '
SetCellText 'E10' '  +6.75%  '
SetCellText 'D11' '5.31'
SetCellText 'E11' '  +11.53%  '
SetCellText 'E12' '  +3.42%  '
SetCellText 'E13' '  +5.57%  '
SetCellText 'D14' '34.23'
SetCellText 'E14' '  +4.12%  '
SetCellText 'E15' '  +1.78%  '
SetCellText 'D16' '3.519.84'
SetCellText 'E16' '  +3.11%  '
SetCellText 'D17' '7.21'
SetCellText 'E17' '  +5.08%  '
SetCellText 'D18' '3.018.10'
SetCellText 'E18' '  +3.12%  '
SetCellText 'D19' '59.869.30'
SetCellText 'E19' '  +3.94%  '
SetCellText 'D20' '439.34'
SetCellText 'E20' '  +5.43%  '
SetCellText 'D21' '13.71'
SetCellText 'E21' '  +3.98%  '
SetCellText 'D22' '0.722'
SetCellText 'E22' '  +5.26%  '
SetCellText 'D23' '7.13'
SetCellText 'E23' '  +2.24%  '
SetCellText 'E24' '  +2.50%  '
SetCellText 'E25' '  +1.32%  '
SetCellText 'E27' '  +12.85%  '
SetCellText 'E28' '  +0.20%  '
SetCellText 'E29' '  +3.46%  '
SetCellText 'D30' '7.85'
SetCellText 'E30' '  +4.88%  '
SetCellText 'B31' 'NEARProtocol'
SetCellText 'C31' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
SetCellText 'D31' '6.30'
SetCellText 'E31' '  +5.13%  '
SetCellText 'B32' 'EthereumClassic'
SetCellText 'C32' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
SetCellText 'D32' '26.07'
SetCellText 'E32' '  +3.76%  '
SetCellText 'D33' '0.104'
SetCellText 'D34' '0.0₃0795'
SetCellText 'E34' '  +16.46%  '
SetCellText 'E35' '  +6.28%  '
SetCellText 'E36' '  +4.96%  '
SetCellText 'D37' '2.13'
SetCellText 'E37' '  +2.84%  '
SetCellText 'D38' '49.22'
SetCellText 'E38' '  +2.32%  '
SetCellText 'D39' '8.66'
SetCellText 'E40' '  +10.63%  '
SetCellText 'D41' '407.03'
SetCellText 'E41' '  +7.37%  '
SetCellText 'E42' '  +2.78%  '
SetCellText 'D43' '2.779.12'
SetCellText 'E43' '  +3.45%  '
SetCellText 'E44' '  -0.44%  '
SetCellText 'D45' '0.254'
SetCellText 'E45' '  +6.71%  '
SetCellText 'E46' '  +0.01%  '
SetCellText 'D47' '122.84'
SetCellText 'E47' '  +0.77%  '
SetCellText 'E48' '  +3.38%  '
SetCellText 'E49' '  +1.54%  '
SetCellText 'D50' '33.98'
SetCellText 'E50' '  +19.91%  '
SetCellText 'D51' '23.69'
SetCellText 'E51' '  +2.76%  '

Write-Output "Applied 80 cell updates to cryptos sheet."
